# Apply the daily update: bump the date and replace each division problem
# in the practice table with its new value. Every old string is unique in
# the document, so a straightforward Find/Replace per pair is safe.

$d = $word.ActiveDocument

$replacements = @(
    @("2026-02-03 Tuesday", "2026-02-04 Wednesday"),
    @("317÷6=", "957÷6="),
    @("985÷6=", "974÷2="),
    @("315÷9=", "829÷2="),
    @("921÷6=", "975÷9="),
    @("336÷2=", "161÷7="),
    @("538÷3=", "780÷8="),
    @("527÷5=", "696÷7="),
    @("615÷2=", "880÷8="),
    @("556÷4=", "846÷2="),
    @("624÷9=", "676÷3="),
    @("592÷4=", "481÷8="),
    @("927÷9=", "611÷4="),
    @("859÷3=", "528÷8="),
    @("830÷4=", "664÷5="),
    @("111÷7=", "768÷7="),
    @("938÷8=", "259÷4="),
    @("318÷6=", "462÷4="),
    @("985÷4=", "786÷9="),
    @("566÷5=", "101÷2="),
    @("863÷5=", "850÷8="),
    @("398÷3=", "642÷2="),
    @("612÷9=", "928÷6="),
    @("411÷5=", "467÷3="),
    @("586÷8=", "203÷3="),
    @("836÷7=", "289÷8=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]

    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, `
                         $true, 1, $false, $new, 2)
}
